$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add the new log entry (Post 48) as row 58 -----------------------------

$row = 58

# Copy the formatting of the last existing data row (57) down onto the new
# row first, so the new cells pick up the same number formats / hyperlink
# style as the rest of the table.
$ws.Range("B57:F57").Copy() | Out-Null
$ws.Range("B58:F58").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Match the shared-string insertion order of the original edit: dev.to link,
# then the title, then the hashnode link (so new sharedStrings indices land
# at 146/147/148 exactly as in the authored workbook).
$ws.Cells.Item($row, 2).Value = 48
$ws.Cells.Item($row, 6).Value = "https://dev.to/rahulmishra05/deadlock-operating-system-m04-p01-5c92"
$ws.Cells.Item($row, 3).Value = "Deadlock | Operating System - M04 P01"
$ws.Cells.Item($row, 4).Value = 44172
$ws.Cells.Item($row, 5).Value = "https://programmingport.hashnode.dev/deadlock-or-operating-system-m04-p01"

# --- Expand the table / autofilter so the new row is included --------------

$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("B10:F58"))

# --- Update the view so the new row is visible / selected -------------------

$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("E58").Select() | Out-Null
